$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new "demo" row to the addresses/cities table: East Palo Alto,
# 2008-2019, location flag "yes", with a note about the Sept 2019 coverage
# change (+0.2 sq miles).
$ws.Range("A36").Value = "East Palo Alto"
$ws.Range("B36").Value = 2008
$ws.Range("C36").Value = 2019
$ws.Range("D36").Value = "yes"
$ws.Range("G36").Value = "change in coverage in Sept 2019. Additional .2 sq miles"

# Reflect the author's final on-screen scroll/selection position after
# adding the row.
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C34").Select() | Out-Null
